# Business logic: allocate parking spots based on employee date of joining.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ssa_employees: correct GECO's "Absent From Date" (date of joining driven update)
$ws1.Range("D4").Value = 42679

# ssa_parking_spots: add "Comments" column describing the allocation/source of each spot
$ws2.Range("D1").Value = "Comments"
$ws2.Range("D1").Font.Bold = $true
$ws2.Range("D2").Value = "NULL"
$ws2.Range("D3").Value = "NULL"
$ws2.Range("D9").Value = "SG Fleet Customer"
$ws2.Range("D10").Value = "PM Partners"
$ws2.Range("D4").Value = "Fleet Partners"
$ws2.Range("D5").Value = "NULL"
$ws2.Range("D6").Value = "NULL"
$ws2.Range("D7").Value = "NULL"
$ws2.Range("D8").Value = "NULL"

# Update the cell selections on both sheets, then make the employees sheet the active tab
$ws2.Range("D8").Select()
$ws1.Activate()
$ws1.Range("D4").Select()
